# Commit: "case with 380 kV done"
# Update voltage-magnitude results (vm_pu) table for Case_3_127 res_bus.
# The slack/reference bus voltage setpoint changed from 1.05 pu to 1.02 pu
# (380 kV case), which in turn changes the recomputed bus voltages for all
# the other buses across every row (time step) of the result table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.063268881355981
$ws.Range("D2").Value = 1.068759701037246
$ws.Range("E2").Value = 1.058427909549006
$ws.Range("F2").Value = 1.077878485871743
$ws.Range("I2").Value = 1.05613503516216
$ws.Range("J2").Value = 1.068235128184039
$ws.Range("K2").Value = 1.071464022579884
$ws.Range("L2").Value = 1.061160171579044
$ws.Range("M2").Value = 1.080558628186115
$ws.Range("N2").Value = 1.069752144594342
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.064672063497058
$ws.Range("D3").Value = 1.069945980059497
$ws.Range("E3").Value = 1.059663834423424
$ws.Range("F3").Value = 1.079297753070762
$ws.Range("I3").Value = 1.056677689077528
$ws.Range("J3").Value = 1.069290912670676
$ws.Range("K3").Value = 1.072465422290839
$ws.Range("L3").Value = 1.062209058128189
$ws.Range("M3").Value = 1.081794205487197
$ws.Range("N3").Value = 1.07080942841605
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.065578216984896
$ws.Range("D4").Value = 1.070712044788009
$ws.Range("E4").Value = 1.060461452159374
$ws.Range("F4").Value = 1.080214986407146
$ws.Range("I4").Value = 1.05702661901967
$ws.Range("J4").Value = 1.069971809872778
$ws.Range("K4").Value = 1.073111283580531
$ws.Range("L4").Value = 1.062885090629076
$ws.Range("M4").Value = 1.082592048752304
$ws.Range("N4").Value = 1.07149129257028
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.065958740406932
$ws.Range("D5").Value = 1.07103373581091
$ws.Range("E5").Value = 1.060796272269235
$ws.Range("F5").Value = 1.080600327775428
$ws.Range("I5").Value = 1.057172784845875
$ws.Range("J5").Value = 1.070257522167068
$ws.Range("K5").Value = 1.073382303325628
$ws.Range("L5").Value = 1.063168662064588
$ws.Range("M5").Value = 1.082927071043878
$ws.Range("N5").Value = 1.071777410608799
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.066022607338557
$ws.Range("D6").Value = 1.071087728043093
$ws.Range("E6").Value = 1.060852460997057
$ws.Range("F6").Value = 1.080665013018323
$ws.Range("I6").Value = 1.057197296076915
$ws.Range("J6").Value = 1.0703054631967
$ws.Range("K6").Value = 1.073427779540873
$ws.Range("L6").Value = 1.063216238029983
$ws.Range("M6").Value = 1.082983300001587
$ws.Range("N6").Value = 1.071825419720192
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.065583303212613
$ws.Range("D7").Value = 1.07071634465724
$ws.Range("E7").Value = 1.060465927992518
$ws.Range("F7").Value = 1.080220136384929
$ws.Range("I7").Value = 1.057028574150886
$ws.Range("J7").Value = 1.069975629677133
$ws.Range("K7").Value = 1.073114906919313
$ws.Range("L7").Value = 1.062888882202945
$ws.Range("M7").Value = 1.082596526865708
$ws.Range("N7").Value = 1.071495117799195
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.063743470223977
$ws.Range("D8").Value = 1.069160931119858
$ws.Range("E8").Value = 1.058846036514902
$ws.Range("F8").Value = 1.078358371166683
$ws.Range("I8").Value = 1.056318885749997
$ws.Range("J8").Value = 1.068592408048358
$ws.Range("K8").Value = 1.071802890538387
$ws.Range("L8").Value = 1.061515203536755
$ws.Range("M8").Value = 1.080976544142652
$ws.Range("N8").Value = 1.070109931837049
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.060487345018541
$ws.Range("D9").Value = 1.066408078065889
$ws.Range("E9").Value = 1.055975148322868
$ws.Range("F9").Value = 1.075068780835847
$ws.Range("I9").Value = 1.055051305166859
$ws.Range("J9").Value = 1.066137388973207
$ws.Range("K9").Value = 1.069474550375283
$ws.Range("L9").Value = 1.059073912482736
$ws.Range("M9").Value = 1.078108965751411
$ws.Range("N9").Value = 1.067651426353217
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.058306655543285
$ws.Range("D10").Value = 1.064564421606643
$ws.Range("E10").Value = 1.054049773400799
$ws.Range("F10").Value = 1.072869317774326
$ws.Range("I10").Value = 1.054194606985337
$ws.Range("J10").Value = 1.06448851893094
$ws.Range("K10").Value = 1.067910973988793
$ws.Range("L10").Value = 1.057432098029287
$ws.Range("M10").Value = 1.076188159766152
$ws.Range("N10").Value = 1.066000214726362
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.057359933264917
$ws.Range("D11").Value = 1.063764023707246
$ws.Range("E11").Value = 1.053213259211293
$ws.Range("F11").Value = 1.071915316682258
$ws.Range("I11").Value = 1.05382083992949
$ws.Range("J11").Value = 1.063771569246983
$ws.Range("K11").Value = 1.067231163565995
$ws.Range("L11").Value = 1.056717701771925
$ws.Range("M11").Value = 1.075354190718128
$ws.Range("N11").Value = 1.065282246891619
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.057007898577012
$ws.Range("D12").Value = 1.063466400700901
$ws.Range("E12").Value = 1.052902109963188
$ws.Range("F12").Value = 1.0715607070568
$ws.Range("I12").Value = 1.05368157988931
$ws.Range("J12").Value = 1.06350480786953
$ws.Range("K12").Value = 1.066978228971191
$ws.Range("L12").Value = 1.05645181309031
$ws.Range("M12").Value = 1.075044072552542
$ws.Range("N12").Value = 1.065015106682411
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.057083428546467
$ws.Range("D13").Value = 1.063530256416616
$ws.Range("E13").Value = 1.052968872146675
$ws.Range("F13").Value = 1.071636783517417
$ws.Range("I13").Value = 1.053711470996485
$ws.Range("J13").Value = 1.063562049762209
$ws.Range("K13").Value = 1.067032503532463
$ws.Range("L13").Value = 1.056508871221319
$ws.Range("M13").Value = 1.075110609678133
$ws.Range("N13").Value = 1.065072429865143
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.05733084176892
$ws.Range("D14").Value = 1.063739428655276
$ws.Range("E14").Value = 1.053187548347711
$ws.Range("F14").Value = 1.071886009709134
$ws.Range("I14").Value = 1.053809337370928
$ws.Range("J14").Value = 1.063749527976986
$ws.Range("K14").Value = 1.067210264574037
$ws.Range("L14").Value = 1.056695734186971
$ws.Range("M14").Value = 1.075328563332609
$ws.Range("N14").Value = 1.06526017432049
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.057483230615915
$ws.Range("D15").Value = 1.063868264024502
$ws.Range("E15").Value = 1.05332222467388
$ws.Range("F15").Value = 1.072039532675192
$ws.Range("I15").Value = 1.053869579463933
$ws.Range("J15").Value = 1.063864979079891
$ws.Range("K15").Value = 1.067319732811413
$ws.Range("L15").Value = 1.056810796163987
$ws.Range("M15").Value = 1.075462805847175
$ws.Range("N15").Value = 1.065375789377202
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.058369433599475
$ws.Range("D16").Value = 1.064617496931386
$ws.Range("E16").Value = 1.054105230134091
$ws.Range("F16").Value = 1.072932596900399
$ws.Range("I16").Value = 1.054219353108522
$ws.Range("J16").Value = 1.064536037186187
$ws.Range("K16").Value = 1.067956031834381
$ws.Range("L16").Value = 1.057479436218953
$ws.Range("M16").Value = 1.076243459575045
$ws.Range("N16").Value = 1.066047800462982
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.058924657949477
$ws.Range("D17").Value = 1.065086908672684
$ws.Range("E17").Value = 1.054595629831205
$ws.Range("F17").Value = 1.073492353468163
$ws.Range("I17").Value = 1.054438001467029
$ws.Range("J17").Value = 1.064956172023233
$ws.Range("K17").Value = 1.068354418920646
$ws.Range("L17").Value = 1.057897919686489
$ws.Range("M17").Value = 1.076732536221218
$ws.Range("N17").Value = 1.066468531939687
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.059248273267731
$ws.Range("D18").Value = 1.065360508127203
$ws.Range("E18").Value = 1.054881400686869
$ws.Range("F18").Value = 1.073818694324066
$ws.Range("I18").Value = 1.05456526439669
$ws.Range("J18").Value = 1.065200942695349
$ws.Range("K18").Value = 1.068586524387224
$ws.Range("L18").Value = 1.058141678763055
$ws.Range("M18").Value = 1.077017590010443
$ws.Range("N18").Value = 1.066713650214232
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.05935857768613
$ws.Range("D19").Value = 1.065453764669305
$ws.Range("E19").Value = 1.054978795434929
$ws.Range("F19").Value = 1.073929941954166
$ws.Range("I19").Value = 1.054608611927397
$ws.Range("J19").Value = 1.065284354719243
$ws.Range("K19").Value = 1.068665621259552
$ws.Range("L19").Value = 1.058224737664698
$ws.Range("M19").Value = 1.077114749509674
$ws.Range("N19").Value = 1.066797180692767
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.058865112277129
$ws.Range("D20").Value = 1.065036566019814
$ws.Range("E20").Value = 1.054543042671459
$ws.Range("F20").Value = 1.073432313045847
$ws.Range("I20").Value = 1.054414570622791
$ws.Range("J20").Value = 1.064911125253819
$ws.Range("K20").Value = 1.068311703417632
$ws.Range("L20").Value = 1.057853055100332
$ws.Range("M20").Value = 1.076680085355936
$ws.Range("N20").Value = 1.066423421198693
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.057257995299831
$ws.Range("D21").Value = 1.063677841520925
$ws.Range("E21").Value = 1.05312316560769
$ws.Range("F21").Value = 1.071812625803975
$ws.Range("I21").Value = 1.053780529963371
$ws.Range("J21").Value = 1.063694332913442
$ws.Range("K21").Value = 1.067157930098641
$ws.Range("L21").Value = 1.056640722393434
$ws.Range("M21").Value = 1.075264390994237
$ws.Range("N21").Value = 1.065204900873624
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.056245334742717
$ws.Range("D22").Value = 1.062821705089925
$ws.Range("E22").Value = 1.052227937109839
$ws.Range("F22").Value = 1.070792807157159
$ws.Range("I22").Value = 1.053379414922031
$ws.Range("J22").Value = 1.062926655120293
$ws.Range("K22").Value = 1.066430058024215
$ws.Range("L22").Value = 1.055875410029737
$ws.Range("M22").Value = 1.074372288458215
$ws.Range("N22").Value = 1.064436132889956
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.056782377004579
$ws.Range("D23").Value = 1.063275737054896
$ws.Range("E23").Value = 1.052702753845492
$ws.Range("F23").Value = 1.071333573075333
$ws.Range("I23").Value = 1.053592288917899
$ws.Range("J23").Value = 1.063333867465484
$ws.Range("K23").Value = 1.066816151050159
$ws.Range("L23").Value = 1.056281410112025
$ws.Range("M23").Value = 1.074845400849266
$ws.Range("N23").Value = 1.064843923523385
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.058892019137886
$ws.Range("D24").Value = 1.065059314317356
$ws.Range("E24").Value = 1.05456680538291
$ws.Range("F24").Value = 1.073459443209586
$ws.Range("I24").Value = 1.054425158850084
$ws.Range("J24").Value = 1.064931480836492
$ws.Range("K24").Value = 1.068331005539915
$ws.Range("L24").Value = 1.057873328511707
$ws.Range("M24").Value = 1.076703786312953
$ws.Range("N24").Value = 1.066443805688628
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.06133084965292
$ws.Range("D25").Value = 1.06712121556863
$ws.Range("E25").Value = 1.056719330186911
$ws.Range("F25").Value = 1.075920317520987
$ws.Range("I25").Value = 1.055381042649139
$ws.Range("J25").Value = 1.066774192858449
$ws.Range("K25").Value = 1.070078458699955
$ws.Range("L25").Value = 1.059707535365554
$ws.Range("M25").Value = 1.078851877693637
$ws.Range("N25").Value = 1.068289134573018
